$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "기술"
$ws.Range("C3").Value = "마이크로소프트, 차세대 AI 칩 공개···‘탈 엔비디아’ 가속"
$ws.Range("D3").Value = "'2026-01-27"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "https://n.news.naver.com/mnews/article/032/0003424054?sid=105"

# Row 4
$ws.Range("B4").Value = "정책"
$ws.Range("C4").Value = "EU, 머스크 AI 그록 '딥페이크 생성' 조사(종합)"
$ws.Range("E4").Value = "https://n.news.naver.com/mnews/article/001/0015867941?sid=104"

# Row 5
$ws.Range("B5").Value = "기업"
$ws.Range("C5").Value = "신한금융, 3천500억원 전략 펀드 조성…AI·에너지·인프라 투자"
$ws.Range("D5").Value = "'2026-01-27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "https://n.news.naver.com/mnews/article/001/0015868879?sid=105"

# Row 6
$ws.Range("B6").Value = "기업"
$ws.Range("C6").Value = "삼성 AI 랩탑 갤럭시북6 시리즈 출시"
$ws.Range("E6").Value = "https://n.news.naver.com/mnews/article/421/0008736850?sid=105"

# Row 7
$ws.Range("B7").Value = "산업"
$ws.Range("C7").Value = "의료AI가 사전에 막았다…뷰노 `"병원내 심정지 46% ↓`""
$ws.Range("E7").Value = "https://n.news.naver.com/mnews/article/003/0013731846?sid=102"

# Row 8
$ws.Range("B8").Value = "정책"
$ws.Range("C8").Value = "과기부, 전북서 피지컬 AI 제조혁신 본격화…지역 AX 사업 추진 논의"
$ws.Range("D8").Value = "'2026-01-26"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "https://n.news.naver.com/mnews/article/022/0004100720?sid=102"

# Row 10
$ws.Range("C10").Value = "[사람과 생각] `"기술봉사가 바꾼 인생`" 한국기술교육대와의 인연으로 취..."
$ws.Range("E10").Value = "https://news.unn.net/news/articleView.html?idxno=589227"

# Row 11
$ws.Range("C11").Value = "작년 사이버 침해 신고 2383건···올해 ‘유출 정보 악용’ 위협 전망"
$ws.Range("E11").Value = "https://n.news.naver.com/mnews/article/032/0003424063?sid=105"

# Row 12
$ws.Range("C12").Value = "과기정통부·KISA, 2025년 사이버 위협 동향 분석 및 2026년 전망 보고서..."
$ws.Range("E12").Value = "http://www.metroseoul.co.kr/article/20260127500352"

# Row 13
$ws.Range("C13").Value = "과기정통부, 피지컬 AI 기반 혁신제품 개발 지원"
$ws.Range("E13").Value = "https://www.koit.co.kr/news/articleView.html?idxno=205131"

# Row 14
$ws.Range("C14").Value = "중소기업 제조 현장에 인공지능 도입 가속화한다"
$ws.Range("E14").Value = "http://www.ikoreanspirit.com/news/articleView.html?idxno=83678"
